# Auto-generated edit script
# Applies profit-recalculation updates to Leve profit columns (H-N) across all 8 sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 7333
$ws.Range("I74").Value = 5999.5
$ws.Range("K74").Value = 5999.5
$ws.Range("M74").Value = -5063.5

$ws.Range("H77").Value = 7333
$ws.Range("I77").Value = 5999.5
$ws.Range("K77").Value = 29997.5
$ws.Range("M77").Value = -25317.5

$ws.Range("H92").Value = 45710.09
$ws.Range("I92").Value = 47791.523
$ws.Range("K92").Value = 47791.523
$ws.Range("M92").Value = -46543.523

$ws.Range("H98").Value = 1693.5555
$ws.Range("I98").Value = 1693.5555
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1693.5555
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -195.5554999999999
$ws.Range("N98").ClearContents()

$ws.Range("H122").Value = 1693.5555
$ws.Range("I122").Value = 1693.5555
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5080.666499999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2630.666499999999
$ws.Range("N122").ClearContents()

$ws.Range("H123").Value = 93987
$ws.Range("J123").Value = 93987
$ws.Range("L123").Value = 93987
$ws.Range("N123").Value = -103787

$ws.Range("H130").Value = 110780
$ws.Range("J130").Value = 110780
$ws.Range("L130").Value = 110780
$ws.Range("N130").Value = -120820

$ws.Range("H131").Value = 2357.0527
$ws.Range("I131").Value = 1477.3572
$ws.Range("J131").Value = 4820.2
$ws.Range("K131").Value = 4432.071599999999
$ws.Range("L131").Value = 14460.6
$ws.Range("M131").Value = 607.9284000000007
$ws.Range("N131").Value = -24540.6

$ws.Range("H138").Value = 5601.755
$ws.Range("I138").Value = 3613.0952
$ws.Range("J138").Value = 6906.8125
$ws.Range("K138").Value = 10839.2856
$ws.Range("L138").Value = 20720.4375
$ws.Range("M138").Value = -5699.285600000001
$ws.Range("N138").Value = -31000.4375

$ws.Range("H141").Value = 11575.846
$ws.Range("I141").Value = 11290.5
$ws.Range("J141").Value = 15000
$ws.Range("K141").Value = 33871.5
$ws.Range("L141").Value = 45000
$ws.Range("M141").Value = -28691.5
$ws.Range("N141").Value = -55360


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 45221.793
$ws.Range("I132").Value = 55154.74
$ws.Range("J132").Value = 7476.6
$ws.Range("K132").Value = 165464.22
$ws.Range("L132").Value = 22429.8
$ws.Range("M132").Value = -162934.22
$ws.Range("N132").Value = -27489.8


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2240.889
$ws.Range("I20").Value = 1752.5714
$ws.Range("K20").Value = 1752.5714
$ws.Range("M20").Value = -1505.5714

$ws.Range("H29").Value = 19166.334
$ws.Range("J29").Value = 18749.5
$ws.Range("L29").Value = 18749.5
$ws.Range("N29").Value = -19327.5

$ws.Range("H75").Value = 7603.25
$ws.Range("I75").Value = 7603.25
$ws.Range("K75").Value = 7603.25
$ws.Range("M75").Value = -6667.25

$ws.Range("H78").Value = 7603.25
$ws.Range("I78").Value = 7603.25
$ws.Range("K78").Value = 22809.75
$ws.Range("M78").Value = -18129.75

$ws.Range("H86").Value = 1866.75
$ws.Range("I86").Value = 1362.75
$ws.Range("K86").Value = 1362.75
$ws.Range("M86").Value = -239.75

$ws.Range("H89").Value = 1866.75
$ws.Range("I89").Value = 1362.75
$ws.Range("K89").Value = 6813.75
$ws.Range("M89").Value = -1197.75

$ws.Range("H94").Value = 8064.421
$ws.Range("I94").Value = 8307.294
$ws.Range("K94").Value = 8307.294
$ws.Range("M94").Value = -7856.294

$ws.Range("H99").Value = 95676
$ws.Range("I99").Value = 129679.5
$ws.Range("J99").Value = 5000
$ws.Range("K99").Value = 129679.5
$ws.Range("L99").Value = 5000
$ws.Range("M99").Value = -128181.5
$ws.Range("N99").Value = -7996

$ws.Range("H105").Value = 4328.875
$ws.Range("I105").Value = 4135.846
$ws.Range("K105").Value = 4135.846
$ws.Range("M105").Value = -2388.846

$ws.Range("H107").Value = 1332.5714
$ws.Range("I107").Value = 1332.5714
$ws.Range("K107").Value = 1332.5714
$ws.Range("M107").Value = 587.4286

$ws.Range("H122").Value = 81302.25
$ws.Range("I122").Value = 709
$ws.Range("J122").Value = 108166.664
$ws.Range("K122").Value = 709
$ws.Range("L122").Value = 108166.664
$ws.Range("M122").Value = 4191
$ws.Range("N122").Value = -117966.664

$ws.Range("H134").Value = 1855.5862
$ws.Range("I134").Value = 1548.7037
$ws.Range("K134").Value = 4646.1111
$ws.Range("M134").Value = -2111.1111


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 77.57895000000001
$ws.Range("J7").Value = 87.71429000000001
$ws.Range("L7").Value = 87.71429000000001
$ws.Range("N7").Value = -313.71429

$ws.Range("H58").Value = 50608.332
$ws.Range("I58").Value = 55650.633
$ws.Range("K58").Value = 55650.633
$ws.Range("M58").Value = -55447.633

$ws.Range("H86").Value = 22207.348
$ws.Range("I86").Value = 37795.5
$ws.Range("K86").Value = 37795.5
$ws.Range("M86").Value = -36672.5

$ws.Range("H88").Value = 17666.666
$ws.Range("J88").Value = 17666.666
$ws.Range("L88").Value = 17666.666
$ws.Range("N88").Value = -18478.666

$ws.Range("H89").Value = 22207.348
$ws.Range("I89").Value = 37795.5
$ws.Range("K89").Value = 188977.5
$ws.Range("M89").Value = -183361.5

$ws.Range("H91").Value = 17666.666
$ws.Range("J91").Value = 17666.666
$ws.Range("L91").Value = 17666.666
$ws.Range("N91").Value = -20474.666

$ws.Range("H123").Value = 81780
$ws.Range("J123").Value = 81780
$ws.Range("L123").Value = 81780
$ws.Range("N123").Value = -91580

$ws.Range("H132").Value = 3628.111
$ws.Range("I132").Value = 3494.8125
$ws.Range("J132").Value = 4694.5
$ws.Range("K132").Value = 10484.4375
$ws.Range("L132").Value = 14083.5
$ws.Range("M132").Value = -7954.4375
$ws.Range("N132").Value = -19143.5

$ws.Range("H134").Value = 113897.664
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws.Range("H136").Value = 50608.332
$ws.Range("I136").Value = 55650.633
$ws.Range("K136").Value = 166951.899
$ws.Range("M136").Value = -164401.899

$ws.Range("H141").Value = 349075.88
$ws.Range("J141").Value = 433383.28
$ws.Range("L141").Value = 433383.28
$ws.Range("N141").Value = -443743.28


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1013.5
$ws.Range("J5").Value = 814.5
$ws.Range("L5").Value = 2443.5
$ws.Range("N5").Value = -2667.5

$ws.Range("H126").Value = 3445
$ws.Range("I126").Value = 3445
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 10335
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -5395
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 1010.6
$ws.Range("I132").Value = 763.25
$ws.Range("K132").Value = 6869.25
$ws.Range("M132").Value = -4339.25

$ws.Range("H135").Value = 1013.5
$ws.Range("J135").Value = 814.5
$ws.Range("L135").Value = 7330.5
$ws.Range("N135").Value = -12400.5

$ws.Range("H140").Value = 2562
$ws.Range("I140").Value = 2613.0908
$ws.Range("K140").Value = 7839.2724
$ws.Range("M140").Value = -2659.2724


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4788.8887
$ws.Range("I70").Value = 4431
$ws.Range("K70").Value = 4431
$ws.Range("M70").Value = -4161

$ws.Range("H73").Value = 4788.8887
$ws.Range("I73").Value = 4431
$ws.Range("K73").Value = 4431
$ws.Range("M73").Value = -3495

$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()

$ws.Range("H122").Value = 3865.3333
$ws.Range("J122").Value = 3904
$ws.Range("L122").Value = 11712
$ws.Range("N122").Value = -16612

$ws.Range("H123").Value = 62960
$ws.Range("J123").Value = 62960
$ws.Range("L123").Value = 62960
$ws.Range("N123").Value = -67860

$ws.Range("H126").Value = 5583.3687
$ws.Range("I126").Value = 4506.0713
$ws.Range("K126").Value = 13518.2139
$ws.Range("M126").Value = -11048.2139

$ws.Range("H132").Value = 67133.375
$ws.Range("I132").Value = 102462.1
$ws.Range("J132").Value = 8252.166999999999
$ws.Range("K132").Value = 307386.3
$ws.Range("L132").Value = 24756.501
$ws.Range("M132").Value = -304856.3
$ws.Range("N132").Value = -29816.501


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()

$ws.Range("H136").Value = 4624.2104
$ws.Range("I136").Value = 3900.5
$ws.Range("J136").Value = 5864.857
$ws.Range("K136").Value = 11701.5
$ws.Range("L136").Value = 17594.571
$ws.Range("M136").Value = -9151.5
$ws.Range("N136").Value = -22694.571


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 100000
$ws.Range("J123").Value = 100000
$ws.Range("L123").Value = 100000
$ws.Range("N123").Value = -109800

$ws.Range("H132").Value = 136333.97
$ws.Range("I132").Value = 147222.92
$ws.Range("K132").Value = 441668.76
$ws.Range("M132").Value = -439138.76

$ws.Range("H136").Value = 5253.625
$ws.Range("J136").Value = 3997
$ws.Range("L136").Value = 11991
$ws.Range("N136").Value = -17091

